$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-20 Friday", "2024-09-21 Saturday"),
    @("639×6=", "341×4="),
    @("301×8=", "864×6="),
    @("822×4=", "227×2="),
    @("923×5=", "935×4="),
    @("924×6=", "106×2="),
    @("131×7=", "257×6="),
    @("928×7=", "376×9="),
    @("508×9=", "221×8="),
    @("334×3=", "451×2="),
    @("708×5=", "823×6="),
    @("141×9=", "918×5="),
    @("434×8=", "453×3="),
    @("780×9=", "843×6="),
    @("513×9=", "371×9="),
    @("581×7=", "231×5="),
    @("593×2=", "534×4="),
    @("311×3=", "685×9="),
    @("735×4=", "396×2="),
    @("584×5=", "172×6="),
    @("270×5=", "644×8="),
    @("828×9=", "151×5="),
    @("118×7=", "732×4="),
    @("468×6=", "573×5="),
    @("737×9=", "823×8="),
    @("706×9=", "994×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
